$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new description value for the data row (row 4, column F)
$ws.Range("F4").Value = "FOOD_DESC_1"

# Update the active selection to F4, matching the edited file's view state
$ws.Activate()
$ws.Range("F4").Select()
